$wb = $excel.ActiveWorkbook

# The tracker moved on to new game dates: rename the last three sheets
# (previously 07-25-24 / 07-26-24 / 07-27-24) to the new dates being tracked.
$wb.Worksheets.Item("07-25-24").Name = "08-04-24"
$wb.Worksheets.Item("07-26-24").Name = "08-06-24"
$wb.Worksheets.Item("07-27-24").Name = "08-07-24"

# Scroll the sheet-tab strip so the newly active sheets are in view.
$excel.ActiveWindow.ScrollWorkbookTabs(13)

# --- 08-04-24: replace the NRFI "Games"/"Score" projections ---
$ws = $wb.Worksheets.Item("08-04-24")
$games = @(
    @("('HOU', 'TB')", 0.751),
    @("('CWS', 'MIN')", 0.746),
    @("('PHI', 'SEA')", 0.742),
    @("('COL', 'SD')", 0.737),
    @("('ATL', 'MIA')", 0.72),
    @("('AZ', 'PIT')", 0.715),
    @("('LAD', 'OAK')", 0.714),
    @("('BOS', 'TEX')", 0.709),
    @("('BAL', 'CLE')", 0.709),
    @("('CHC', 'STL')", 0.675),
    @("('DET', 'KC')", 0.67),
    @("('NYY', 'TOR')", 0.67),
    @("('CIN', 'SF')", 0.407),
    @("('MIL', 'WSH')", 0.189),
    @("('LAA', 'NYM')", 0.117)
)
for ($i = 0; $i -lt $games.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $games[$i][0]
    $ws.Cells.Item($row, 2).Value = $games[$i][1]
}

# --- 08-06-24: replace the NRFI "Games"/"Score" projections ---
$ws = $wb.Worksheets.Item("08-06-24")
$games = @(
    @("('CWS', 'OAK')", 0.747),
    @("('LAD', 'PHI')", 0.738),
    @("('SF', 'WSH')", 0.73),
    @("('PIT', 'SD')", 0.72),
    @("('COL', 'NYM')", 0.719),
    @("('CHC', 'MIN')", 0.711),
    @("('HOU', 'TEX')", 0.705),
    @("('CIN', 'MIA')", 0.6840000000000001),
    @("('ATL', 'MIL')", 0.5610000000000001),
    @("('AZ', 'CLE')", 0.535),
    @("('BOS', 'KC')", 0.521),
    @("('BAL', 'TOR')", 0.478),
    @("('DET', 'SEA')", 0.038),
    @("('LAA', 'NYY')", -0.008),
    @("('STL', 'TB')", -0.126)
)
for ($i = 0; $i -lt $games.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $games[$i][0]
    $ws.Cells.Item($row, 2).Value = $games[$i][1]
}

# --- 08-07-24: replace the NRFI "Games"/"Score" projections (shorter list) ---
$ws = $wb.Worksheets.Item("08-07-24")
$games = @(
    @("('COL', 'NYM')", 0.727),
    @("('DET', 'SEA')", 0.708),
    @("('LAD', 'PHI')", 0.705),
    @("('STL', 'TB')", 0.702),
    @("('PIT', 'SD')", 0.6840000000000001),
    @("('SF', 'WSH')", 0.653),
    @("('CIN', 'MIA')", 0.635),
    @("('BOS', 'KC')", 0.598),
    @("('AZ', 'CLE')", 0.5590000000000001),
    @("('HOU', 'TEX')", 0.485),
    @("('CHC', 'MIN')", 0.463),
    @("('BAL', 'TOR')", 0.319),
    @("('ATL', 'MIL')", 0.163),
    @("('CWS', 'OAK')", -0.08699999999999999)
)
for ($i = 0; $i -lt $games.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $games[$i][0]
    $ws.Cells.Item($row, 2).Value = $games[$i][1]
}
# This sheet went from 16 rows to 15 rows - clear the now-stale trailing rows.
$ws.Range("A16:B17").ClearContents()
